$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 46

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "08/01/2025"
$ws.Cells.Item($row, 2).Value = 556.875
$ws.Cells.Item($row, 3).Value = 0.08978675645342311
$ws.Cells.Item($row, 4).Value = 50
